$d = $word.ActiveDocument

$pairs = @(
    @("96×42=", "96×60="),
    @("87×56=", "41×28="),
    @("61×22=", "15×99="),
    @("28×90=", "70×52="),
    @("35×11=", "21×19="),
    @("59×43=", "97×16="),
    @("95×53=", "60×87="),
    @("28×49=", "87×41="),
    @("17×92=", "86×73="),
    @("74×61=", "87×49="),
    @("84×32=", "58×55="),
    @("56×17=", "65×55="),
    @("73×73=", "55×67="),
    @("48×21=", "62×55="),
    @("31×31=", "90×73="),
    @("75×30=", "77×57="),
    @("35×70=", "63×99="),
    @("27×50=", "62×48="),
    @("65×78=", "58×51="),
    @("46×62=", "46×93="),
    @("80×68=", "98×27="),
    @("52×52=", "69×98="),
    @("68×25=", "55×48="),
    @("13×76=", "11×61="),
    @("65×49=", "60×99=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
